$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text revisions ("Wrap up Semaster 1") - update dialogue/script lines in column B.
# Order matches the authored edit sequence so the shared-strings table is rebuilt
# in the same order as the source commit.
$ws.Range("B21").Value = "Alright…"
$ws.Range("B39").Value = "I deflected a bullet."
$ws.Range("B44").Value = "I am sure it was not a illusion. I just can feel a deeper connection with everything."
$ws.Range("B38").Value = "Is everything ok?"
$ws.Range("B40").Value = "Hah, you must be kidding."
$ws.Range("B43").Value = "Just focus on the mission, Lyu. It was weird and we will figure that out later."
$ws.Range("B24").Value = "Search around. The key should be nearby."

# Reflect the author's final cursor position/selection in the sheet view.
$ws.Range("H19").Select() | Out-Null
